# Generate Report for Handback
#
# This localization status report is regenerated once the de-de handback
# has been verified to be in sync with en-US. The "Status" columns move
# from "Ready for handoff" to "Handed back: in sync with en-US", the
# handback timestamps are refreshed, and the stale "handback file is not
# the latest" error details are cleared now that everything is current.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet: status, refreshed handback datetime, cleared error ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-09-07 01:07:30"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet: status, refreshed handback datetime, cleared error ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-09-07 01:07:38"
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments (status columns widened, error column narrowed) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.17   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 29.17   # column F (de-de status)

$wsZhCn.Columns.Item(3).ColumnWidth = 29.17        # column C (Status)
$wsZhCn.Columns.Item(16).ColumnWidth = 12.83       # column P (Error Detail)

$wsDeDe.Columns.Item(3).ColumnWidth = 29.17        # column C (Status)
$wsDeDe.Columns.Item(16).ColumnWidth = 12.83       # column P (Error Detail)
